# Univariate Model 1 - slide edits
#   1. Remove the "Slope " run from the Slope oval (id=47), leaving an
#      empty paragraph.
#   2. Remove the "Intercept" run from the Intercept oval (id=155),
#      leaving the trailing space run untouched.
#   3. Delete the floating "I-Mean" / "S-Mean" textboxes (id=200 / id=202).
#   4. Nudge the two "Construct Name" textboxes (id=10 / id=62) to their
#      new positions.
#   5. Add four new centered textboxes: "Slope", "M = S-Mean",
#      "M = I-Mean", "Intercept" at their new dedicated positions.

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# EMU<->point helper. PowerPoint's Left/Top/Width/Height properties are
# expressed in points (1 pt = 12700 EMU) and are backed by 32-bit floats,
# so a plain round-trip through $emu/12700 is frequently off the target
# integer EMU by a hair. A small nudge keeps it on the correct side of
# the rounding boundary for the (modest) magnitudes used on this slide.
$EMU_PER_POINT = 12700
$EMU_EPS = 0.3
function ToPt($emu) {
    return ($emu + $EMU_EPS) / $EMU_PER_POINT
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) "Slope " oval -> strip the run, keep the paragraph (and its
#    endParaRPr) intact.
$slopeOval = Get-ShapeById $s.Shapes 47
$slopeOval.TextFrame.TextRange.Paragraphs(3, 1).Text = ""

# 2) "Intercept" oval -> strip only the "Intercept" run, keep the " "
#    run that follows it.
$interceptOval = Get-ShapeById $s.Shapes 155
$interceptOval.TextFrame.TextRange.Characters(3, 9).Text = ""

# 3) Remove the two stray mean labels.
(Get-ShapeById $s.Shapes 200).Delete()
(Get-ShapeById $s.Shapes 202).Delete()

# 4) Reposition the two "Construct Name" textboxes.
$constructName1 = Get-ShapeById $s.Shapes 10
$constructName1.Left = ToPt 751916
$constructName1.Top = ToPt 2983776

$constructName2 = Get-ShapeById $s.Shapes 62
$constructName2.Left = ToPt 862820
$constructName2.Top = ToPt 837041

# 5) Add the four new labels, cloned from an existing "Construct Name"
#    textbox so they pick up the same run/paragraph/body formatting
#    (centered alignment, spAutoFit, etc.).
$template = Get-ShapeById $s.Shapes 62

function Add-Label($template, $text, $name, $x, $y) {
    $clone = $template.Duplicate().Item(1)
    $clone.Name = $name
    $clone.TextFrame.TextRange.Text = $text
    $clone.Left = ToPt $x
    $clone.Top = ToPt $y
}

Add-Label $template "Slope" "TextBox 62" 825170 1082147
Add-Label $template "M = S-Mean" "TextBox 63" 825170 1324530
Add-Label $template "M = I-Mean" "TextBox 64" 763306 3447580
Add-Label $template "Intercept" "TextBox 65" 795912 3206749
